$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 147
$ws1.Range("F4").Value = 36
$ws1.Range("F7").Value = 38
$ws1.Range("F9").Value = 557
$ws1.Range("F13").Value = 162
$ws1.Range("F14").Value = 25
$ws1.Range("F17").Value = 109
$ws1.Range("F18").Value = 5057
$ws1.Range("F19").Value = 55
$ws1.Range("F20").Value = 834
$ws1.Range("F21").Value = 115
$ws1.Range("F22").Value = 2262
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 2106

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 147
$ws4.Range("F4").Value = 36
$ws4.Range("F7").Value = 38
$ws4.Range("F9").Value = 557
$ws4.Range("F13").Value = 162
$ws4.Range("F14").Value = 25
$ws4.Range("F17").Value = 109
$ws4.Range("F18").Value = 5057
$ws4.Range("F20").Value = 55
$ws4.Range("F22").Value = 834
$ws4.Range("F23").Value = 115
$ws4.Range("F24").Value = 2262
$ws4.Range("F26").Value = 27
$ws4.Range("F27").Value = 2106
